# CRM-4377 : Add Inactive SF also in SF Document List Download with Active and Inactive status
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before column B ("Address") to hold the vendor company name.
$ws.Columns.Item(2).Insert()

# Copy the header style used by the other "vendor:xxx" header cells (e.g. District, GST*)
# onto the newly inserted "Company Name" header cell.
$ws.Range("E1").Copy() | Out-Null
$ws.Range("B1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("B1").Value = "Company Name"
$ws.Range("B2").Value = "{vendor:company_name}"

# Match column width roughly to column A
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(1).ColumnWidth

# Append a new "Status" column at the end (now column AG after the insert)
$lastCol = $ws.Cells(1, $ws.Columns.Count).End(-4159).Column  # xlToLeft
$newCol = $lastCol + 1

$ws.Cells(1, $newCol).Value = "Status"
$ws.Cells(2, $newCol).Value = "{vendor:active_status}"

# Give the new header/value cells the same style used by the other header/value columns.
$ws.Range("E1").Copy() | Out-Null
$ws.Cells(1, $newCol).PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("A2").Copy() | Out-Null
$ws.Cells(2, $newCol).PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Application.CutCopyMode = 0

# Leave the selection where the author left it when saving.
$ws.Range("B5").Select() | Out-Null
